# "BISAGRAS 1842 DISMAY" price list update:
#  - bump the list date in A1 by one month (2024-04-24 -> 2024-05-24)
#  - refresh the PRECIO column for the five "Bisagra 1842" SKUs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45436

$ws.Range("D33").Value = 277.464
$ws.Range("D34").Value = 307.393
$ws.Range("D35").Value = 355.55
$ws.Range("D39").Value = 379.403
$ws.Range("D40").Value = 570.456
